$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new yuv_* columns BG:BX
$headers = @(
  "yuv_y_mean", "yuv_y_median", "yuv_y_variance", "yuv_y_std_dev", "yuv_y_25", "yuv_y_75", "yuv_u_mean", "yuv_u_median", "yuv_u_variance", "yuv_u_std_dev", "yuv_u_25", "yuv_u_75", "yuv_v_mean", "yuv_v_median", "yuv_v_variance", "yuv_v_std_dev", "yuv_v_25", "yuv_v_75"
)
$headerCols = @("BG", "BH", "BI", "BJ", "BK", "BL", "BM", "BN", "BO", "BP", "BQ", "BR", "BS", "BT", "BU", "BV", "BW", "BX")

# Copy the existing header style (from A1, which already carries the
# bold/centered/bordered header format, style index 1) onto the new
# header cells before writing their text, so the new headers share the
# exact same style as the rest of row 1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("BG1:BX1").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $headerCols.Length; $i++) {
  $ws.Range($headerCols[$i] + "1").Value = $headers[$i]
}

# Data rows 2:22 for columns BG:BX
$data = @{
  "2" = @(51.02854625651144, 39, 1614.383248020209, 40.17938834801009, 31, 51, 126.7295789718628, 128, 18.18695630187145, 4.264616782534095, 128, 128, 130.3387759526571, 128, 52.08011642569901, 7.216655487530149, 128, 128)
  "3" = @(60.21131947057391, 51, 1283.713154333453, 35.8289429698038, 41, 65, 121.9487761903479, 128, 72.92286701188964, 8.539488685623375, 116, 128, 138.9517624065284, 128, 177.5055480445853, 13.32312080725028, 128, 150)
  "4" = @(54.23931422607637, 44, 1476.896198701306, 38.43040721487746, 33, 59, 127.0047152837118, 128, 15.65922770760736, 3.957174207386801, 128, 128, 129.8872728347778, 128, 43.95832536464248, 6.630107492691388, 128, 128)
  "5" = @(63.21710892787601, 57, 1186.066728516082, 34.43931951296486, 45, 68, 121.8383954366048, 128, 130.1741873606792, 11.40939031502907, 122, 128, 138.164667447408, 128, 326.8593086963999, 18.07925077807153, 128, 146)
  "6" = @(65.26228724510544, 58, 1291.895184130998, 35.94294345391037, 45, 72, 126.0012976328532, 128, 43.13220173615247, 6.56751107621087, 128, 128, 131.3423665364583, 128, 108.8130766484482, 10.43135066271134, 128, 128)
  "7" = @(81.60848817211372, 71, 1363.594742408876, 36.92688373541527, 59, 93, 122.1528600056966, 128, 138.0444998261455, 11.74923401018745, 128, 128, 137.1279433568319, 128, 316.2135171882184, 17.78239346061768, 128, 128)
  "8" = @(68.25345021468003, 62, 807.5120681069233, 28.41675681894264, 52, 77, 119.640360635123, 128, 160.5535997257944, 12.67097469517615, 107, 128, 140.4784920998451, 128, 337.0554489888107, 18.35906993801186, 128, 160)
  "9" = @(73.07390639056339, 65, 1146.6291508048, 33.8619129820629, 52, 85, 118.5313874246679, 128, 166.7278265936989, 12.91231298388088, 106, 128, 142.346434142897, 128, 352.4103937496778, 18.77259688348093, 128, 162)
  "10" = @(84.66241152422101, 82, 1322.70041725695, 36.3689485310883, 59, 105, 114.0180567087372, 128, 246.7167629583285, 15.70722009008369, 99, 128, 149.6535512151613, 128, 523.2349146780023, 22.87432872628183, 128, 173)
  "11" = @(83.11840290907979, 78, 1043.306075021197, 32.30024883837889, 64, 93, 125.7879486083984, 128, 59.87909693620168, 7.738158497743612, 128, 128, 131.2087167104085, 128, 121.1546248449686, 11.00702615809414, 128, 128)
  "12" = @(95.20793541473456, 93, 1003.184422443058, 31.67308672111163, 73, 113, 124.996148109436, 128, 67.3523471599072, 8.206847577474996, 128, 128, 132.4037745793661, 128, 139.0231814897541, 11.79080919571486, 128, 128)
  "13" = @(115.6294208667128, 111, 1187.698709428388, 34.46300493904134, 92, 133, 122.5339393615723, 128, 176.3684388790134, 13.28037796446372, 128, 128, 135.1483821868896, 128, 288.6540494431239, 16.9898219367692, 128, 128)
  "14" = @(93.64590091439908, 91, 1023.565905953612, 31.99321656154022, 71, 113, 118.300250814487, 128, 182.5719701244901, 13.51191955735713, 104, 128, 141.3945496691889, 128, 337.6246138997089, 18.37456431863648, 128, 162)
  "15" = @(136.5870167398562, 134, 1197.620420979426, 34.6066528427617, 113, 160, 115.5106943242674, 128, 280.9928149361689, 16.76284030038373, 99, 128, 143.9788966581698, 128, 434.8087008335709, 20.85206706380859, 128, 166)
  "16" = @(148.0959565386791, 150, 1211.781181280123, 34.81064752744659, 125, 174, 116.9730920940347, 128, 258.0455616084456, 16.06379661252114, 103, 128, 142.360382150104, 128, 411.3512929644501, 20.28179708419474, 128, 162)
  "17" = @(142.0877163442713, 145, 1195.751370765808, 34.5796380947778, 121, 166, 114.2799731215093, 128, 292.0745383874604, 17.0901883660614, 98, 128, 145.8338931588577, 128, 457.846700082011, 21.39735264190435, 128, 167)
  "18" = @(138.5684528969394, 140, 1135.227874020394, 33.69314283382293, 118, 160, 124.133929570516, 128, 117.9385791360107, 10.85995299879381, 128, 128, 133.27090771993, 128, 206.6474876083563, 14.3752386974393, 128, 128)
  "19" = @(142.5104241846531, 141, 965.5847520975916, 31.07385962666356, 124, 161, 113.9332907096061, 128, 281.54417880735, 16.77927825644923, 99, 128, 146.5224611255821, 128, 448.3209170780729, 21.17359008477478, 128, 167)
  "20" = @(144.04752615442, 143, 1042.645388841806, 32.2900199572841, 124, 164, 122.7012424468994, 128, 158.3390886510259, 12.58328608317501, 128, 128, 134.6793622970581, 128, 243.1977810249169, 15.59479980714459, 128, 128)
  "21" = @(138.2128079629797, 137, 1004.159835642087, 31.68848111920303, 119, 153, 119.376017430845, 128, 190.4071254805602, 13.79880884281539, 109, 128, 139.7894189718328, 128, 329.8148055611053, 18.16080410006961, 128, 157)
  "22" = @(164.4903306435105, 165, 636.2851579068138, 25.22469341551674, 148, 182, 117.4711913442789, 128, 196.7827020386028, 14.02792579245424, 105, 128, 141.6723067692769, 128, 308.1596903313408, 17.55447778577707, 128, 160)
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + $r).Value = $vals[$i]
  }
}
